# Strip the stray trailing "16" that was appended to every Bible reference
# in column A (e.g. "Galatians 1:116" -> "Galatians 1:1") so the references
# are human readable again. Column B (verse text) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value.EndsWith("16")) {
        $cell.Value2 = $value.Substring(0, $value.Length - 2)
    }
}
